# Change order of Low-Cal and Cal in the Protocol sheet:
# Old layout: Cal MSB | Cal LSB | Low Cal MSB | Low Cal LSB | High Cal MSB | High Cal LSB
# New layout: Low Cal (MSB) | Low Cal (LSB) |  Cal (MSB) | Cal (LSB) | High Cal (MSB) | High Cal (LSB)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# Header row 8 (Master->Slave table)
$ws.Range("F8").Value = "Low Cal (MSB)"
$ws.Range("G8").Value = "Low Cal (LSB)"
$ws.Range("H8").Value = " Cal (MSB)"
$ws.Range("I8").Value = "Cal (LSB)"
$ws.Range("J8").Value = "High Cal (MSB)"
$ws.Range("K8").Value = "High Cal (LSB)"

# Header row 13 (Slave->Master table)
$ws.Range("F13").Value = "Low Cal (MSB)"
$ws.Range("G13").Value = "Low Cal (LSB)"
$ws.Range("H13").Value = " Cal (MSB)"
$ws.Range("I13").Value = "Cal (LSB)"
$ws.Range("J13").Value = "High Cal (MSB)"
$ws.Range("K13").Value = "High Cal (LSB)"

$ws.Range("M9").Select()
